$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.195.99'
$ws.Range('E2').Value = '  -0.64%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.859.63'
$ws.Range('E3').Value = '  -1.26%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7038'
$ws.Range('E5').Value = '  -1.18%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.19'
$ws.Range('E6').Value = '  -0.07%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07771'
$ws.Range('E9').Value = '  -3.29%  '

# Row 10
$ws.Range('E10').Value = '  -4.47%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08004'
$ws.Range('E11').Value = '  -3.93%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.851.71'
$ws.Range('E12').Value = '  -1.98%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.172'
$ws.Range('E13').Value = '  -1.40%  '

# Row 14
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.32'
$ws.Range('E14').Value = '  +0.38%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6947'
$ws.Range('E15').Value = '  -3.60%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.363'
$ws.Range('E16').Value = '  +0.29%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008275'
$ws.Range('E17').Value = '  -2.17%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.173.79'
$ws.Range('E18').Value = '  -0.76%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '250.76'
$ws.Range('E19').Value = '  +3.89%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.119.96'
$ws.Range('E20').Value = '  -1.46%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.08'
$ws.Range('E21').Value = '  -1.34%  '

# Row 22
$ws.Range('E22').Value = '  -0.02%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.495'
$ws.Range('E23').Value = '  -4.48%  '

# Row 24
$ws.Range('E24').Value = '  -0.02%  '

# Row 25
$ws.Range('E25').Value = '  -2.27%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.961'
$ws.Range('E26').Value = '  -0.94%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.48'
$ws.Range('E27').Value = '  -2.78%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.73'
$ws.Range('E28').Value = '  +0.86%  '

# Row 29
$ws.Range('E29').Value = '  -0.65%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.278'
$ws.Range('E30').Value = '  -3.15%  '

# Row 31
$ws.Range('E31').Value = '  -2.03%  '

# Row 32
$ws.Range('E32').Value = '  +0.88%  '

# Row 33
$ws.Range('E33').Value = '  -2.28%  '

# Row 34
$ws.Range('E34').Value = '  -3.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7421'
$ws.Range('E35').Value = '  -1.09%  '

# Row 36
$ws.Range('E36').Value = '  -2.54%  '

# Row 37
$ws.Range('E37').Value = '  -0.05%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01862'
$ws.Range('E38').Value = '  -1.40%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.251.00'
$ws.Range('E39').Value = '  -3.19%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.732'
$ws.Range('E40').Value = '  -0.53%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.240'
$ws.Range('E41').Value = '  -5.52%  '

# Row 42
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.83'
$ws.Range('E42').Value = '  -0.98%  '

# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8950'
$ws.Range('E43').Value = '  -2.79%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.28'
$ws.Range('E44').Value = '  -3.80%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.04%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000129'
$ws.Range('E46').Value = '  +0.33%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.020.43'
$ws.Range('E47').Value = '  -1.34%  '

# Row 48
$ws.Range('E48').Value = '  -0.40%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.776'
$ws.Range('E49').Value = '  -1.77%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.395'
$ws.Range('E50').Value = '  -1.10%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4285'
$ws.Range('E51').Value = '  -2.59%  '
